$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Initial Values"
# ---------------------------------------------------------------------------
$wsIV = $wb.Worksheets.Item("Initial Values")

# Utilization: 0.982 -> 0.91
$wsIV.Range("C6").Value2 = 0.91

# Wage/basic consumption: replace the fixed 11.9 with a live ratio formula,
# and bump the display precision to three decimals (matches style index 12).
$wsIV.Range("C11").NumberFormat = "0.000"
$wsIV.Range("C11").Formula = "=6522/1322"

# ---------------------------------------------------------------------------
# Sheet "GDP" -- drop parameter "q" entirely and retune several coefficients
# ---------------------------------------------------------------------------
$wsGDP = $wb.Worksheets.Item("GDP")

$rowVals = @{
    5  = 0.73
    6  = 4.72
    7  = 0.68
    8  = 0.57999999999999996
    9  = 0.20599999999999999
    10 = 0.02
    13 = 0.49
}
foreach ($r in $rowVals.Keys) {
    $wsGDP.Range("C" + $r + ":L" + $r).Value2 = $rowVals[$r]
}

# Remove the "q" row entirely: clear the label and all of its values.
$wsGDP.Range("B11").ClearContents() | Out-Null
$wsGDP.Range("C11:L11").ClearContents() | Out-Null

# ---------------------------------------------------------------------------
# Sheet "Exports" -- refresh scenario figures, columns D:M
# ---------------------------------------------------------------------------
$wsExp = $wb.Worksheets.Item("Exports")

$exportRows = @{
    4  = @(0.06, 0.05, 0.03, 0.01, 0.01, 0.03, 0.03, 0.04, 0.05, 0.05)
    5  = @(0.06, 0.04, 0.02, 0,    0,    0.02, 0.03, 0.04, 0.04, 0.04)
    6  = @(0.06, 0.06, 0.06, 0.05, 0.05, 0.04, 0.04, 0.04, 0.04, 0.04)
    7  = @(0.06, 0.05, -0.04, 0,   0.04, 0.05, 0.06, 0.06, 0.06, 0.06)
    8  = @(0.13200000000000001, 0.13, 0.12, 0.1,  0.09, 0.09, 0.09, 0.09, 0.09, 0.09)
    9  = @(0.13200000000000001, 0.1,  0.07, 0.07, 0.07, 0.07, 0.07, 0.07, 0.07, 0.07)
    10 = @(0.13200000000000001, 0.13, 0.15, 0.17, 0.19, 0.2,  0.2,  0.2,  0.2,  0.2)
    11 = @(0.13200000000000001, 0.13, 0.25, 0.2,  0.15, 0.15, 0.15, 0.15, 0.15, 0.15)
}
$cols = @("D", "E", "F", "G", "H", "I", "J", "K", "L", "M")
foreach ($r in $exportRows.Keys) {
    $vals = $exportRows[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $wsExp.Range($cols[$i] + $r).Value2 = $vals[$i]
    }
}

# ---------------------------------------------------------------------------
# Sheet "Investment and Saving" -- retune a handful of coefficients
# ---------------------------------------------------------------------------
$wsIS = $wb.Worksheets.Item("Investment and Saving")

$isRows = @{
    5  = 0.058
    7  = 0.367
    8  = 0.1
    15 = 0.393
}
foreach ($r in $isRows.Keys) {
    $wsIS.Range("C" + $r + ":L" + $r).Value2 = $isRows[$r]
}

# ---------------------------------------------------------------------------
# Restore cursor / selection positions on each touched sheet. The final
# selection must land on "Exports" so it remains the active tab.
# ---------------------------------------------------------------------------
$wsIV.Range("G10").Select() | Out-Null
$wsGDP.Range("G22").Select() | Out-Null
$wsIS.Range("E19").Select() | Out-Null
$wsExp.Range("G17").Select() | Out-Null
